$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the rubric text in C14 (add "adecuados" to the description)
$ws.Range("C14").Value = "Uso de mensajes de error y/o validación adecuados"

# Swap the point values between D9 and D14
$ws.Range("D9").Value = 0.5
$ws.Range("D14").Value = 0.75

# Update the selected cell shown when the file was last saved
$ws.Range("D10").Select()

# Update the saved window position/size
$excel.ActiveWindow.Left = -24120
$excel.ActiveWindow.Top = -120
$excel.ActiveWindow.Width = 24240
$excel.ActiveWindow.Height = 13140
